# LCV icin OTV, hurda ve kredi eklendi
# - Rename "esneklik" -> "esneklik (binek)"
# - Insert a new worksheet "LCV esneklik" right after it, with its own
#   elasticity value (mirrors the "esneklik" layout: Degisken/Deger header
#   plus one data row).
# - Leave "segments" and "kleit_segment_elasticity" untouched (they just
#   shift right in the tab order as a side effect of the insertion).

$wb = $excel.ActiveWorkbook

$binek = $wb.Worksheets.Item("esneklik")
$binek.Name = "esneklik (binek)"

$lcv = $wb.Worksheets.Add($null, $binek)
$lcv.Name = "LCV esneklik"

$lcv.Range("A1").Value = "Degisken"
$lcv.Range("B1").Value = "Deger"
$lcv.Range("A1:B1").Font.Bold = $true

$lcv.Range("A2").Value = "LCV  fiyat-talep esnekligi"
$lcv.Range("B2").Value = -2.8

$lcv.Columns.Item(1).ColumnWidth = 16.92

# Restore the view state each tab had in the authored workbook: the
# previously-active "esneklik" sheet keeps a plain range selection, while
# the newly inserted sheet becomes the active tab.
$binek.Range("A1:B4").Select()
$lcv.Range("G22").Select()
$lcv.Activate()
